$d = $word.ActiveDocument

# 1) Insert a new changelog entry (Compact style, numId 1001) right before
#    the existing "2019-06-13: Revised grade bundles ..." entry, i.e. right
#    after the "Until the census date ..." paragraph. Inserting a paragraph
#    break right before that entry makes the new (empty) paragraph inherit
#    its "Compact" / numbered-list formatting, matching the target diff.
$rng = $d.Content
$found = $rng.Find.Execute("2019-06-13: Revised grade bundles", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $targetPara = $rng.Paragraphs.Item(1)
    $targetIndex = $targetPara.Index
    $insertPt = $d.Range($rng.Start, $rng.Start)
    $insertPt.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($targetIndex)
    $newPara.Range.Text = "2019-06-26: Increased the quota for LTA re-takes per week (from 3 to 5) and clarified Mastery homework extension."
}

# 2) Bump the weekly LTA re-take quota from 3 to 5.
$d.Content.Find.Execute(
    "no more than 3 LTA re-takes per week",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "no more than 5 LTA re-takes per week", 2) | Out-Null

# 3) Clarify that the Mastery assignment extension token cost is per module.
$d.Content.Find.Execute(
    "24-hour deadline extension on Mastery assignment (1 token/assignment);",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "24-hour deadline extension on Mastery assignment (1 token/assignment module);", 2) | Out-Null
